$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale hyperlink on M3 before the row shifts beneath it.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$M$3') {
        $h.Delete()
    }
}

# Row 3 (the second data record) is being removed entirely.
$ws.Rows(3).Delete()

# The remaining record's Reference Code (O2) changes value.
$ws.Range("O2").Value = "Email45"

# Selection moves from D3 to D2 to track the now-shorter sheet.
$ws.Range("D2").Select()
